$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching original inline-string formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.854.43"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "2.248.62"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "306.74"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "95.95"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "0.573"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "35.27"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "7.23"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "2.592.26"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.344.75"
$ws.Range("E15").Value = "  +4.66%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.842"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "13.64"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "44.606.15"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "0.0₃0950"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").Value = "12.00"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "6.28"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "65.48"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "239.77"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "1.98"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").Value = "  +3.84%  "
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Value = "37.76"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "151.00"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").Value = "0.0799"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E35").Value = "  -7.90%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  +4.31%  "
$ws.Range("D39").Value = "15.08"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "3.80"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").Value = "0.0302"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "1.856.06"
$ws.Range("E44").Value = "  +7.67%  "
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").Value = "  +14.04%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.191"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "79.82"
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("D48").Value = "99.24"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "69.33"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "54.68"
$ws.Range("E51").Value = "  -0.03%  "

# Restore default (unstyled) cell style now that text values are locked in
$ws.Range("D2:D51").Style = "Normal"
